$d = $word.ActiveDocument
$shp = $d.InlineShapes.Item(1)

# --- Crop the picture (right edge ~12.831% of the native 96dpi width) ---
$shp.PictureFormat.CropRight = 80.8353

# --- Resize to the final displayed extent (468pt x 118.8827559055118pt) ---
# Temporarily unlock the aspect ratio so Width/Height can be set
# independently, then restore the lock so the saved XML keeps
# noChangeAspect="1" on the graphicFrameLocks.
$shp.LockAspectRatio = $false
$shp.Width = 468
$shp.Height = 118.8827559055118
$shp.LockAspectRatio = $true

# --- Remove the picture outline/border ---
$shp.Line.Visible = $false

# --- Alt text (Word's auto-generated picture description) ---
# Use literal "&#10;" numeric character references (rather than real
# newline bytes) so the saved attribute round-trips as two real
# newlines instead of being XML-whitespace-normalized to spaces.
$shp.AlternativeText = "A blue screen with white text&#10;&#10;Description automatically generated"
